# Generate Report for Archive
#
# The localization status text "Ready for handoff" is renamed to
# "In Translation" everywhere it appears (Overview!E2:F3 and the
# "Status" column of the per-language detail tables on the zh-cn and
# de-de sheets). The status columns are also narrowed to fit the new,
# shorter label.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $current = $cell.Value()
        if ($oldStatus -eq $current) {
            $cell.Value = $newStatus
        }
    }
}

# Shrink the status columns to match the shorter text (best fit available
# through the column-width API): Overview columns E & F ("zh-cn"/"de-de"
# status) and column C ("Status") on the zh-cn / de-de detail sheets.
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth
